# Flip the sign (+ -> -) of the "material recycled" figures for every
# component (Nd, Dy, Cu, Si columns) on every year sheet.
#
# Each yearly worksheet has the same A1:E5 layout:
#   row2 -> Nd, row4 -> Cu, row5 -> Si, columns B..E -> Generator Onshore /
#   Generator Offshore / Panel / Wires. Only C2, B4, C4, E4 and D5 ever hold
#   a non-zero figure; negating an already-zero cell is a harmless no-op,
#   so the same unconditional negate is applied to every sheet.

$wb = $excel.ActiveWorkbook
$cellRefs = @("C2", "B4", "C4", "E4", "D5")

$sheetCount = $wb.Worksheets.Count
for ($i = 1; $i -le $sheetCount; $i++) {
    $ws = $wb.Worksheets.Item($i)
    foreach ($ref in $cellRefs) {
        $cell = $ws.Range($ref)
        $current = $cell.Value2
        $cell.Value = (0 - $current)
    }
}
